# "Generate Report for Handback" - refresh the handoff/handback timestamps
# for the 0c44ee9e... file (row 2) on the zh-cn and de-de language sheets,
# and roll the Overview sheet's "Latest HO Xliff Generate Date" up to the
# newest of those per-language handoff timestamps.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-16 04:45:56"
$zhcn.Range("K2").Value = "2016-08-16 04:46:26"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-16 04:46:03"
$dede.Range("K2").Value = "2016-08-16 04:46:33"

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-16 04:46:03"
